# Apply grading corrections to the Kanagala Lab Exam 03 grading workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20 (addProduct() method): revise points + add grading comment ---
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1) For fetching the customers and checking the customer is null it" + [char]0x2019 + "s a wrong way. You have to fetch products of that customer and check that products are null or not."

# --- Row 34 (CustomerMappingTest Class): revise points + replace grading comment ---
$ws.Range("F34").WrapText = $true
$ws.Range("F34").Value = "(-3)I have changed your addProduct() code and run the test cases then 3 test cases failed but I didn" + [char]0x2019 + "t deducted any points for remaining test cases" + [char]10
$ws.Rows.Item(34).RowHeight = 30
$ws.Range("E34").Value = 4

# --- Column F: widen to fit the long grading comments ---
$ws.Columns.Item(6).ColumnWidth = 150.6

# --- Restore editor's scroll/selection state ---
$excel.Goto($ws.Range("D16"))
$ws.Range("H33").Select()

$wb.Save()
